$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "93.735.38"
$ws.Range("E2").Value = "  -1.44%  "
$ws.Range("D3").Value = "3.323.78"
$ws.Range("E3").Value = "  -3.30%  "
$ws.Range("E4").Value = "  +0.03%  "
$s = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.13"
$ws.Range("D5").Style = $s
$ws.Range("E5").Value = "  -3.70%  "
$s = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "618.51"
$ws.Range("D6").Style = $s
$s = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.39"
$ws.Range("D7").Style = $s
$ws.Range("E7").Value = "  -2.98%  "
$s = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.388"
$ws.Range("D8").Style = $s
$ws.Range("E8").Value = "  -3.98%  "
$ws.Range("E9").Value = "  -0.07%  "
$s = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.933"
$ws.Range("D10").Style = $s
$ws.Range("E10").Value = "  -5.95%  "
$ws.Range("D11").Value = "3.323.57"
$ws.Range("E11").Value = "  -3.18%  "
$s = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.03"
$ws.Range("D12").Style = $s
$ws.Range("E12").Value = "  +0.32%  "
$s = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.193"
$ws.Range("D13").Style = $s
$ws.Range("E13").Value = "  -2.19%  "
$ws.Range("B14").Value = "Toncoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$s = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.96"
$ws.Range("D14").Style = $s
$ws.Range("E14").Value = "  -1.87%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "93.520.54"
$ws.Range("E15").Value = "  -1.37%  "
$ws.Range("D16").Value = "3.953.13"
$ws.Range("E16").Value = "  -3.23%  "
$ws.Range("E17").Value = "  -4.57%  "
$s = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.10"
$ws.Range("D18").Style = $s
$ws.Range("E18").Value = "  -3.64%  "
$ws.Range("D19").Value = "3.322.49"
$ws.Range("E19").Value = "  -3.49%  "
$s = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.19"
$ws.Range("D20").Style = $s
$ws.Range("E20").Value = "  -4.38%  "
$s = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.90"
$ws.Range("D21").Style = $s
$ws.Range("E21").Value = "  -5.71%  "
$s = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.47"
$ws.Range("D22").Style = $s
$ws.Range("E22").Value = "  +8.82%  "
$s = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "494.66"
$ws.Range("D23").Style = $s
$ws.Range("E23").Value = "  -1.68%  "
$s = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.454"
$ws.Range("D24").Style = $s
$ws.Range("E24").Value = "  -10.46%  "
$ws.Range("E25").Value = "  -4.83%  "
$s = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.15"
$ws.Range("D26").Style = $s
$ws.Range("E26").Value = "  -6.58%  "
$s = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "91.36"
$ws.Range("D27").Style = $s
$ws.Range("E27").Value = "  -0.03%  "
$s = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.69"
$ws.Range("D28").Style = $s
$ws.Range("E28").Value = "  -3.16%  "
$ws.Range("D29").Value = "3.507.82"
$ws.Range("E29").Value = "  -3.07%  "
$ws.Range("E30").Value = "  +0.22%  "
$s = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.05"
$ws.Range("D31").Style = $s
$ws.Range("E31").Value = "  -5.71%  "
$s = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.138"
$ws.Range("D32").Style = $s
$ws.Range("E32").Value = "  +0.71%  "
$s = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.63"
$ws.Range("D33").Style = $s
$ws.Range("E33").Value = "  -4.14%  "
$s = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = $s
$ws.Range("E34").Value = "  +0.48%  "
$ws.Range("E35").Value = "  -5.05%  "
$s = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "28.30"
$ws.Range("D36").Style = $s
$ws.Range("E36").Value = "  -8.20%  "
$s = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.530"
$ws.Range("D37").Style = $s
$ws.Range("E37").Value = "  -6.25%  "
$s = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "528.98"
$ws.Range("D38").Style = $s
$ws.Range("E38").Value = "  +3.47%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$s = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.41"
$ws.Range("D39").Style = $s
$ws.Range("E39").Value = "  -3.74%  "
$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$s = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = $s
$ws.Range("E40").Value = "  +0.06%  "
$s = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.148"
$ws.Range("D41").Style = $s
$ws.Range("E41").Value = "  -1.38%  "
$s = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.36"
$ws.Range("D42").Style = $s
$ws.Range("E42").Value = "  -4.99%  "
$s = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.866"
$ws.Range("D43").Style = $s
$ws.Range("E43").Value = "  -4.43%  "
$ws.Range("B44").Value = "MantraDAO"
$ws.Range("C44").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$s = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.78"
$ws.Range("D44").Style = $s
$ws.Range("E44").Value = "  +6.92%  "
$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$s = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "24.06"
$ws.Range("D45").Style = $s
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$s = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0417"
$ws.Range("D46").Style = $s
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("B47").Value = "ImmutableX"
$ws.Range("C47").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$s = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.68"
$ws.Range("D47").Style = $s
$ws.Range("E47").Value = "  -0.74%  "
$s = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.39"
$ws.Range("D48").Style = $s
$ws.Range("E48").Value = "  -2.25%  "
$s = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "52.92"
$ws.Range("D49").Style = $s
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("E50").Value = "  -2.16%  "
$s = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.97"
$ws.Range("D51").Style = $s
$ws.Range("E51").Value = "  -0.04%  "
